$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Đơn sale chính" -- update service-group labels (G column) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G2").Value = "Tiểu phẫu"
$ws1.Range("G3").Value = "Tiêm"
$ws1.Range("G4").Value = "Tiêm"
$ws1.Range("G5").Value = "Phun xăm"

# --- Sheet 2: "Đơn thu nợ" -- rebuild with new column layout (A1:AB4) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Clear()

$ws2.Range("A1").Value = "Ngày thực hiện"
$ws2.Range("B1").Value = "Ngày thu"
$ws2.Range("C1").Value = "notion id"
$ws2.Range("D1").Value = "Tiền tố"
$ws2.Range("E1").Value = "Mã đơn thu nợ"
$ws2.Range("F1").Value = "Cơ sở"
$ws2.Range("G1").Value = "id đơn nợ"
$ws2.Range("H1").Value = "Lượng thu"
$ws2.Range("I1").Value = "Đơn nợ"
$ws2.Range("J1").Value = "Nguồn khách"
$ws2.Range("K1").Value = "Sale chính"
$ws2.Range("L1").Value = "Đơn giá gốc"
$ws2.Range("M1").Value = "Sale phụ"
$ws2.Range("N1").Value = "Upsale"
$ws2.Range("O1").Value = "Bác sĩ 1"
$ws2.Range("P1").Value = "Bác sĩ 2"
$ws2.Range("Q1").Value = "Thanh toán lần đầu"
$ws2.Range("R1").Value = "Đã thanh toán"
$ws2.Range("S1").Value = "Tỉ lệ chiết khấu sale chính"
$ws2.Range("T1").Value = "Tỉ lệ chiết khấu sale phụ"
$ws2.Range("U1").Value = "id sale chính"
$ws2.Range("V1").Value = "id sale phụ"
$ws2.Range("W1").Value = "id bác sĩ 1"
$ws2.Range("X1").Value = "id bác sĩ 2"
$ws2.Range("Y1").Value = "Chiết khấu bác sĩ 1"
$ws2.Range("Z1").Value = "Chiết khấu bác sĩ 2"
$ws2.Range("AA1").Value = "Chiết khấu sale chính"
$ws2.Range("AB1").Value = "Chiết khấu sale phụ"
$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "04-04-2024"
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "07-05-2024"
$ws2.Range("C2").Value = "dfafe21b-d5c8-4925-8327-15d40efe6033"
$ws2.Range("D2").Value = "TN"
$ws2.Range("E2").Value = 137
$ws2.Range("F2").Value = "SÓC TRĂNG"
$ws2.Range("G2").Value = "375900ba-870d-4372-be18-f3821db1f765"
$ws2.Range("H2").Value = 10000000
$ws2.Range("I2").Value = "HD-LUXURY-262"
$ws2.Range("J2").Value = "Cá nhân"
$ws2.Range("K2").Value = "Thạch Hoàng Nhân"
$ws2.Range("L2").Value = 10000000
$ws2.Range("O2").Value = "CTV Ngoài"
$ws2.Range("R2").Value = 10000000
$ws2.Range("S2").Value = 0.1
$ws2.Range("T2").Value = 0
$ws2.Range("U2").Value = "cca1354d-d585-4e09-8845-dc6dadbcb631"
$ws2.Range("W2").Value = "7bb857c9-f973-440b-88f2-97e138ee6082"
$ws2.Range("Y2").Value = 1000000
$ws2.Range("Z2").Value = 0
$ws2.Range("AA2").Value = 1000000
$ws2.Range("AB2").Value = 0
$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = "04-29-2024"
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "07-05-2024"
$ws2.Range("C3").Value = "24a4f213-618b-4ea3-b604-5b10881d56f7"
$ws2.Range("D3").Value = "TN"
$ws2.Range("E3").Value = 139
$ws2.Range("F3").Value = "SÓC TRĂNG"
$ws2.Range("G3").Value = "803377a1-0326-4663-9ed5-acaf19520399"
$ws2.Range("H3").Value = 2000000
$ws2.Range("I3").Value = "HD-LUXURY-356"
$ws2.Range("J3").Value = "CTV"
$ws2.Range("K3").Value = "Thạch Hoàng Nhân"
$ws2.Range("L3").Value = 35000000
$ws2.Range("O3").Value = "Phạm Thanh Hoàng"
$ws2.Range("Q3").Value = 15000000
$ws2.Range("R3").Value = 24000000
$ws2.Range("S3").Value = 0
$ws2.Range("T3").Value = 0
$ws2.Range("U3").Value = "cca1354d-d585-4e09-8845-dc6dadbcb631"
$ws2.Range("W3").Value = "a73ea60d-3de1-4e9b-aa7b-f22fda5742bd"
$ws2.Range("Y3").Value = 200000
$ws2.Range("Z3").Value = 0
$ws2.Range("AA3").Value = 0
$ws2.Range("AB3").Value = 0
$ws2.Range("D4").Value = "Tổng"
$ws2.Range("E4").Value = 2
$ws2.Range("H4").Value = 12000000
$ws2.Range("L4").Value = 45000000
$ws2.Range("N4").Value = 0
$ws2.Range("Q4").Value = 15000000
$ws2.Range("R4").Value = 34000000
$ws2.Range("S4").Value = 0.1
$ws2.Range("T4").Value = 0
$ws2.Range("Y4").Value = 1200000
$ws2.Range("Z4").Value = 0
$ws2.Range("AA4").Value = 1000000
$ws2.Range("AB4").Value = 0
